# Assignment 1: Hybrid Sort - expand the explanation paragraph.
#
# The original document has a single body paragraph (paragraph 3) that
# explains the hybrid sort. This edit:
#   1. Splits that paragraph after "...O(nlogn). " and inserts a new
#      sentence about why insertion sort was chosen.
#   2. Adds a blank paragraph.
#   3. Starts a new paragraph ("What I found with this is that even ...")
#      that reuses the old "with the running time..." text, adds a new
#      "(depending on the runsize value)" aside (keeping the _GoBack
#      bookmark positioned right after "value"), and rewrites the final
#      sentence about why the hybrid sort is faster and mentions runsize.
#   4. Adds another blank paragraph and a final paragraph holding a single
#      tab character.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 0. The existing "_GoBack" bookmark sits at the very end of paragraph 3.
#    Remove it now; we re-add it later at its new target location once
#    the surrounding text exists.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 1. "This hybrid sort uses both insertion sort and merge sort, but even"
#    is removed and replaced by the new "For this project, ..." sentence.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("This hybrid sort uses both insertion sort and merge sort, but even")
$rng.Text = "For this project, I decided to use insertion sort to make a hybrid sort because amongst the other in-place sorts we went over, insertion is the easiest to understand and I feel has the best " + [char]0x201C + "best case" + [char]0x201D + " running time (being equal with bubble sort but I just prefer the logic of insertion). "

# ---------------------------------------------------------------------
# 2. Break the paragraph right after the sentence just inserted, add an
#    empty paragraph, and break again so "What I found..." starts a new
#    paragraph of its own.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("logic of insertion). ")
$rng.Collapse(0)
$rng.InsertParagraphAfter()

$rng = $d.Content
$rng.Find.Execute("logic of insertion). ")
$rng.Collapse(0)
$rng.MoveStart(1, 1) | Out-Null
$rng.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 3. Prefix the remaining "with the running time..." text with
#    "What I found with this is that even".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" with the running time of insertion sort being O(n")
$rng.Collapse(1)
$rng.InsertBefore("What I found with this is that even")

# ---------------------------------------------------------------------
# 4. "insertion is faster than merge sort" (the clause right after
#    "for a certain size, ") becomes "different sorting algorithms
#    could be faster than merge sort".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("insertion is faster than merge sort")
$rng.Text = "different sorting algorithms could be faster than merge sort"

# ---------------------------------------------------------------------
# 5. "really optimize the running time. " becomes "really optimize the
#    running time with the perfect runsize value!".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("really optimize the running time. ")
$rng.Text = "really optimize the running time with the perfect runsize value!"

# ---------------------------------------------------------------------
# 6. Insert the "(depending on the runsize value)" aside right after
#    "...faster than merge sort" (before the "! It is because..."
#    sentence). The _GoBack bookmark belongs right after "value" and
#    before the closing ")".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("faster than merge sort! It is because for a certain size, different sorting algorithms could be faster than merge sort")
$rng.Collapse(1)
$rng.MoveEnd(1, "faster than merge sort".Length) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" (depending on the runsize value)")

$rng = $d.Content
$rng.Find.Execute(" (depending on the runsize value)")
$rng.Collapse(0)
$rng.MoveStart(1, -1) | Out-Null
$rng.MoveEnd(1, -1) | Out-Null
$d.Bookmarks.Add("_GoBack", $rng)

# ---------------------------------------------------------------------
# 7. Trailing blank paragraph + a final paragraph holding a tab.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Collapse(0)
$rng.MoveStart(1, -1) | Out-Null
$rng.InsertParagraphAfter()

$rng = $d.Content
$rng.Collapse(0)
$rng.MoveStart(1, -1) | Out-Null
$rng.InsertParagraphAfter()

$rng = $d.Content
$rng.Collapse(0)
$rng.MoveStart(1, -1) | Out-Null
$rng.InsertAfter([char]9)

Write-Output "done"
